$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "tabernacle" -> "dwelling place" (first occurrence only - the one
#    wrapped in the commentRangeStart/End for Brett Slote's "Dwelling
#    place?" comment), then drop a fresh "_GoBack" bookmark right
#    before it (mirrors Word re-anchoring its last-edit marker there).
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("tabernacle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Text = "dwelling place"

$r1again = $d.Content
$r1again.Find.Execute("dwelling place", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($r1again.Start, $r1again.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 2) Old "_GoBack" bookmark sitting between "Allelui" and "a." gets
#    folded away: merge the two runs back into a single "Alleluia."
#    run (removing the stray bookmark in the process).
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Allelui", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $r2.End
$straddle = $d.Range($splitPos - 1, $splitPos + 1)
$straddle.Delete()
$fix = $d.Range($splitPos - 1, $splitPos - 1)
$fix.InsertAfter("ia")

Write-Output "done"
